# Update the vocabulary table: 13 new German/English word pairs inserted in
# alphabetical order by the German word (column A), the translation for
# "furchtbar" corrected from "abscheulich" to "terrible", and the alphabet
# reference letters in column C re-aligned to their new rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A/B (and column C letter markers where present) for all 50 rows ---
$ws.Cells.Item(1, 1).Value = "abenteuer"
$ws.Cells.Item(1, 2).Value = "adventure"
$ws.Cells.Item(1, 3).Value = "A"

$ws.Cells.Item(2, 1).Value = "anhalten"
$ws.Cells.Item(2, 2).Value = "to stop"

$ws.Cells.Item(3, 1).Value = "ankreuzen"
$ws.Cells.Item(3, 2).Value = "to check"

$ws.Cells.Item(4, 1).Value = "ansehen"
$ws.Cells.Item(4, 2).Value = "to watch"

$ws.Cells.Item(5, 1).Value = "befangen"
$ws.Cells.Item(5, 2).Value = "biased"
$ws.Cells.Item(5, 3).Value = "B"

$ws.Cells.Item(6, 1).Value = "beklagen"
$ws.Cells.Item(6, 2).Value = "to mourn"

$ws.Cells.Item(7, 1).Value = "beschreiben"
$ws.Cells.Item(7, 2).Value = "to describe"

$ws.Cells.Item(8, 1).Value = "besonders"
$ws.Cells.Item(8, 2).Value = "especially"

$ws.Cells.Item(9, 1).Value = "chemie"
$ws.Cells.Item(9, 2).Value = "chemical"
$ws.Cells.Item(9, 3).Value = "C"

$ws.Cells.Item(10, 1).Value = "dusche"
$ws.Cells.Item(10, 2).Value = "shower"
$ws.Cells.Item(10, 3).Value = "D"

$ws.Cells.Item(11, 1).Value = "eigenschaften"
$ws.Cells.Item(11, 2).Value = "characteristics"
$ws.Cells.Item(11, 3).Value = "E"

$ws.Cells.Item(12, 1).Value = "einladen"
$ws.Cells.Item(12, 2).Value = "to invite"

$ws.Cells.Item(13, 1).Value = "erhalten"
$ws.Cells.Item(13, 2).Value = "to receive"

$ws.Cells.Item(14, 1).Value = "Entscheidung"
$ws.Cells.Item(14, 2).Value = "decision"

$ws.Cells.Item(15, 1).Value = "euch"
$ws.Cells.Item(15, 2).Value = "you"

$ws.Cells.Item(16, 1).Value = "familie"
$ws.Cells.Item(16, 2).Value = "family"
$ws.Cells.Item(16, 3).Value = "F"

$ws.Cells.Item(17, 1).Value = "feiern sie"
$ws.Cells.Item(17, 2).Value = "to celebrate"

$ws.Cells.Item(18, 1).Value = "furchtbar"
$ws.Cells.Item(18, 2).Value = "terrible"

$ws.Cells.Item(19, 1).Value = "geben"
$ws.Cells.Item(19, 2).Value = "to give"
$ws.Cells.Item(19, 3).Value = "G"

$ws.Cells.Item(20, 1).Value = "halten"
$ws.Cells.Item(20, 2).Value = "to hold"
$ws.Cells.Item(20, 3).Value = "H"

$ws.Cells.Item(21, 1).Value = "idee"
$ws.Cells.Item(21, 2).Value = "idea"
$ws.Cells.Item(21, 3).Value = "I"

$ws.Cells.Item(22, 1).Value = "jacke"
$ws.Cells.Item(22, 2).Value = "jacket"
$ws.Cells.Item(22, 3).Value = "J"

$ws.Cells.Item(23, 1).Value = "Katze"
$ws.Cells.Item(23, 2).Value = "cat"
$ws.Cells.Item(23, 3).Value = "K"

$ws.Cells.Item(24, 1).Value = "komisch"
$ws.Cells.Item(24, 2).Value = "funny"

$ws.Cells.Item(25, 1).Value = "langsam"
$ws.Cells.Item(25, 2).Value = "slow"
$ws.Cells.Item(25, 3).Value = "L"

$ws.Cells.Item(26, 1).Value = "letzte"
$ws.Cells.Item(26, 2).Value = "last"

$ws.Cells.Item(27, 1).Value = "Meldung"
$ws.Cells.Item(27, 2).Value = "message"
$ws.Cells.Item(27, 3).Value = "M"

$ws.Cells.Item(28, 1).Value = "meinen"
$ws.Cells.Item(28, 2).Value = "to mean"

$ws.Cells.Item(29, 1).Value = "möglich"
$ws.Cells.Item(29, 2).Value = "possible"

$ws.Cells.Item(30, 1).Value = "nacht"
$ws.Cells.Item(30, 2).Value = "night"
$ws.Cells.Item(30, 3).Value = "N"

$ws.Cells.Item(31, 1).Value = "opa"
$ws.Cells.Item(31, 2).Value = "grandpa"
$ws.Cells.Item(31, 3).Value = "O"

$ws.Cells.Item(32, 1).Value = "paar"
$ws.Cells.Item(32, 2).Value = "pair"
$ws.Cells.Item(32, 3).Value = "P"

$ws.Cells.Item(33, 1).Value = "qualität"
$ws.Cells.Item(33, 2).Value = "quality"
$ws.Cells.Item(33, 3).Value = "Q"

$ws.Cells.Item(34, 1).Value = "rabatt"
$ws.Cells.Item(34, 2).Value = "discount"
$ws.Cells.Item(34, 3).Value = "R"

$ws.Cells.Item(35, 1).Value = "ruhig"
$ws.Cells.Item(35, 2).Value = "calm"

$ws.Cells.Item(36, 1).Value = "sache"
$ws.Cells.Item(36, 2).Value = "thing"
$ws.Cells.Item(36, 3).Value = "S"

$ws.Cells.Item(37, 1).Value = "sich verabschieden"
$ws.Cells.Item(37, 2).Value = "to farewell"

$ws.Cells.Item(38, 1).Value = "tasche"
$ws.Cells.Item(38, 2).Value = "bag"
$ws.Cells.Item(38, 3).Value = "T"

$ws.Cells.Item(39, 1).Value = "und"
$ws.Cells.Item(39, 2).Value = "and"
$ws.Cells.Item(39, 3).Value = "U"

$ws.Cells.Item(40, 1).Value = "unterstützen"
$ws.Cells.Item(40, 2).Value = "to support"

$ws.Cells.Item(41, 1).Value = "verrückt"
$ws.Cells.Item(41, 2).Value = "crazy"
$ws.Cells.Item(41, 3).Value = "V"

$ws.Cells.Item(42, 1).Value = "viele"
$ws.Cells.Item(42, 2).Value = "many"

$ws.Cells.Item(43, 1).Value = "wem"
$ws.Cells.Item(43, 2).Value = "who"
$ws.Cells.Item(43, 3).Value = "W"

$ws.Cells.Item(44, 1).Value = "wer"
$ws.Cells.Item(44, 2).Value = "who"

$ws.Cells.Item(45, 1).Value = "würst"
$ws.Cells.Item(45, 2).Value = "sausage"

$ws.Cells.Item(46, 1).Value = "zeit"
$ws.Cells.Item(46, 2).Value = "time"
$ws.Cells.Item(46, 3).Value = "Z"

$ws.Cells.Item(47, 1).Value = "zwielichtig"
$ws.Cells.Item(47, 2).Value = "dodgy"

$ws.Cells.Item(48, 1).Value = "ähnlich"
$ws.Cells.Item(48, 2).Value = "similar"
$ws.Cells.Item(48, 3).Value = "Ä"

$ws.Cells.Item(49, 1).Value = "öffnen"
$ws.Cells.Item(49, 2).Value = "to open"
$ws.Cells.Item(49, 3).Value = "Ö"

$ws.Cells.Item(50, 1).Value = "über"
$ws.Cells.Item(50, 2).Value = "about"
$ws.Cells.Item(50, 3).Value = "Ü"

# --- Clear stale column C letter markers from rows that no longer carry one ---
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(6, 3).ClearContents()
$ws.Cells.Item(7, 3).ClearContents()
$ws.Cells.Item(12, 3).ClearContents()
$ws.Cells.Item(13, 3).ClearContents()
$ws.Cells.Item(14, 3).ClearContents()
$ws.Cells.Item(15, 3).ClearContents()
$ws.Cells.Item(17, 3).ClearContents()
$ws.Cells.Item(24, 3).ClearContents()
$ws.Cells.Item(28, 3).ClearContents()
$ws.Cells.Item(29, 3).ClearContents()
$ws.Cells.Item(35, 3).ClearContents()
$ws.Cells.Item(37, 3).ClearContents()

# --- The bold "final letter" formatting used to sit on the old last two rows
#     (36/37); those are now ordinary rows, so their special formatting is
#     removed, and applied instead to the new final two rows (49/50), copied
#     from the untouched bold-styled cell D27.
$ws.Range("D27").Copy()
$ws.Range("C49").PasteSpecial(-4122)
$ws.Range("C50").PasteSpecial(-4122)
$ws.Range("C36").ClearFormats()
$ws.Range("C37").ClearFormats()
$excel.CutCopyMode = $false

# --- Cursor / selection ends on B21, as in the edited workbook ---
$ws.Range("B21").Select()
